$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: use data from before-row 10
$ws.Cells.Item(2, 4).Value = 44630
$ws.Cells.Item(2, 10).Value = 90
$ws.Cells.Item(2, 11).Value = 2500
$ws.Cells.Item(2, 12).Value = 3000
$ws.Cells.Item(2, 13).Value = 2722
$ws.Cells.Item(2, 15).Value = "Región Metropolitana"
$ws.Cells.Item(2, 16).Value = 454

# Row 3: use data from before-row 7
$ws.Cells.Item(3, 4).Value = 44644
$ws.Cells.Item(3, 10).Value = 140
$ws.Cells.Item(3, 11).Value = 2500
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = 2786
$ws.Cells.Item(3, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(3, 16).Value = 464

# Row 4: use data from before-row 2
$ws.Cells.Item(4, 4).Value = 44659
$ws.Cells.Item(4, 10).Value = 90
$ws.Cells.Item(4, 11).Value = 2500
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 13).Value = 2722
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 454

# Row 6: use data from before-row 4
$ws.Cells.Item(6, 4).Value = 44631
$ws.Cells.Item(6, 10).Value = 110
$ws.Cells.Item(6, 11).Value = 3000
$ws.Cells.Item(6, 12).Value = 3500
$ws.Cells.Item(6, 13).Value = 3273
$ws.Cells.Item(6, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(6, 16).Value = 546

# Row 7: use data from before-row 6
$ws.Cells.Item(7, 4).Value = 44658
$ws.Cells.Item(7, 10).Value = 180
$ws.Cells.Item(7, 11).Value = 2500
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = 2778
$ws.Cells.Item(7, 15).Value = "Región Metropolitana"
$ws.Cells.Item(7, 16).Value = 463

# Row 9: use data from before-row 3
$ws.Cells.Item(9, 4).Value = 44643
$ws.Cells.Item(9, 10).Value = 90
$ws.Cells.Item(9, 11).Value = 2800
$ws.Cells.Item(9, 12).Value = 3000
$ws.Cells.Item(9, 13).Value = 2911
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 485

# Row 10: use data from before-row 11
$ws.Cells.Item(10, 4).Value = 44672
$ws.Cells.Item(10, 10).Value = 140
$ws.Cells.Item(10, 11).Value = 3000
$ws.Cells.Item(10, 12).Value = 3500
$ws.Cells.Item(10, 13).Value = 3286
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 548

# Row 11: use data from before-row 12
$ws.Cells.Item(11, 4).Value = 44685
$ws.Cells.Item(11, 10).Value = 150
$ws.Cells.Item(11, 11).Value = 3000
$ws.Cells.Item(11, 12).Value = 3500
$ws.Cells.Item(11, 13).Value = 3267
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 544

# Row 12: use data from before-row 9
$ws.Cells.Item(12, 4).Value = 44650
$ws.Cells.Item(12, 10).Value = 130
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 12).Value = 3500
$ws.Cells.Item(12, 13).Value = 3308
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 551

